## VisitFreqResults.xlsx edit
## - GroupA Weekday: selection changes to A1:D1 (no data change)
## - GroupA Weekend: fills in B/C/D data for rows 2-16, drops tabSelected,
##   sets selection to B17
## - adds three new sheets: "GroupB Weekday", "GroupB Weekend", "Sheet3"
##   with their own data / selections, and activates "Sheet3" last.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: GroupA Weekday -- selection only change
# ---------------------------------------------------------------------------
$wkday = $wb.Worksheets.Item(1)
$wkday.Activate()
$wkday.Range("A1:D1").Select()

# ---------------------------------------------------------------------------
# Sheet 2: GroupA Weekend -- fill rows 2-16 with B/C/D values
# ---------------------------------------------------------------------------
$wkend = $wb.Worksheets.Item(2)
$wkend.Activate()

$wkendRows = @(
    @{r=2;  b="21 28 29 "; c=1.1305270000000001;  d=0.001802},
    @{r=3;  b="21 29 ";    c=0.88156000000000001; d=0.01704},
    @{r=4;  b="21 29 ";    c=1.0830869999999999;  d=0.006338},
    @{r=5;  b="NULL";      c="NULL";               d="NULL"},
    @{r=6;  b="21 29 ";    c=0.69742999999999999; d=0.01583},
    @{r=7;  b="21 29";     c=0.73680000000000001; d=0.0114},
    @{r=8;  b="21 28 29 "; c=1.1997370000000001;  d=0.003036},
    @{r=9;  b="21 29 ";    c=0.61319999999999997; d=0.9456},
    @{r=10; b="NULL";      c="NULL";               d="NULL"},
    @{r=11; b="21 29";     c=0.86224299999999998; d=0.005636},
    @{r=12; b="21 29";     c=0.88705999999999996; d=0.01007},
    @{r=13; b="21 29";     c=0.45472000000000001; d=0.02016},
    @{r=14; b=29;          c=1.31172;              d=0.01111},
    @{r=15; b="21 29";     c=0.80020000000000002; d=0.0159},
    @{r=16; b="21 29";     c=0.27417000000000002; d=0.05277}
)

foreach ($row in $wkendRows) {
    $r = $row.r
    $wkend.Cells.Item($r, 2).Value = $row.b
    $wkend.Cells.Item($r, 3).Value = $row.c
    $wkend.Cells.Item($r, 4).Value = $row.d
}

$wkend.Range("B17").Select()

# ---------------------------------------------------------------------------
# Sheet 3 (new): GroupB Weekday
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$grpBWkday = $wb.Worksheets.Add($null, $last)
$grpBWkday.Name = "GroupB Weekday"
$grpBWkday.Activate()

# header row
$grpBWkday.Cells.Item(1, 1).Value = "Job Type ID"
$grpBWkday.Cells.Item(1, 2).Value = "Facility Numbers"
$grpBWkday.Cells.Item(1, 3).Value = "Estimate"
$grpBWkday.Cells.Item(1, 4).Value = "Std Error"

for ($i = 1; $i -le 33; $i++) {
    $grpBWkday.Cells.Item($i + 1, 1).Value = $i
}

$grpBWkdayData = @(
    @{r=2;  b=6; c=0.63712500000000005; d=0.003553},
    @{r=15; b=6; c=0.73249900000000001; d=0.004632},
    @{r=21; b=6; c=0.28042;              d=0.01079},
    @{r=31; b=6; c=0.50468800000000003; d=0.008602}
)
$grpBWkdayDataRows = @{}
foreach ($item in $grpBWkdayData) { $grpBWkdayDataRows[$item.r] = $item }

for ($r = 2; $r -le 34; $r++) {
    if ($grpBWkdayDataRows.ContainsKey($r)) {
        $item = $grpBWkdayDataRows[$r]
        $grpBWkday.Cells.Item($r, 2).Value = $item.b
        $grpBWkday.Cells.Item($r, 3).Value = $item.c
        $grpBWkday.Cells.Item($r, 4).Value = $item.d
    } else {
        $grpBWkday.Cells.Item($r, 2).Value = "NULL"
        $grpBWkday.Cells.Item($r, 3).Value = "NULL"
        $grpBWkday.Cells.Item($r, 4).Value = "NULL"
    }
}

$grpBWkday.Range("A1:D1").Select()

# ---------------------------------------------------------------------------
# Sheet 4 (new): GroupB Weekend
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$grpBWkend = $wb.Worksheets.Add($null, $last)
$grpBWkend.Name = "GroupB Weekend"
$grpBWkend.Activate()

$grpBWkend.Cells.Item(1, 1).Value = "Job Type ID"
$grpBWkend.Cells.Item(1, 2).Value = "Facility Numbers"
$grpBWkend.Cells.Item(1, 3).Value = "Estimate"
$grpBWkend.Cells.Item(1, 4).Value = "Std Error"

for ($i = 1; $i -le 33; $i++) {
    $grpBWkend.Cells.Item($i + 1, 1).Value = $i
}

# rows 2-20 get B = NULL unless overridden below; rows 21+ have no B/C/D
# (row 21 has only B, no C/D) except the two data rows at 2 and 15.
for ($r = 2; $r -le 20; $r++) {
    $grpBWkend.Cells.Item($r, 2).Value = "NULL"
    $grpBWkend.Cells.Item($r, 3).Value = "NULL"
    $grpBWkend.Cells.Item($r, 4).Value = "NULL"
}

$grpBWkend.Cells.Item(2, 2).Value = 6
$grpBWkend.Cells.Item(2, 3).Value = -1.2346999999999999
$grpBWkend.Cells.Item(2, 4).Value = 0.1251

$grpBWkend.Cells.Item(15, 2).Value = 6
$grpBWkend.Cells.Item(15, 3).Value = -1.5764
$grpBWkend.Cells.Item(15, 4).Value = 0.2386

$grpBWkend.Cells.Item(21, 2).Value = 6

$grpBWkend.Application.Goto($grpBWkend.Range("A8"))
$grpBWkend.Range("B22:D34").Select()

# ---------------------------------------------------------------------------
# Sheet 5 (new): Sheet3 -- empty, becomes the active sheet
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet3 = $wb.Worksheets.Add($null, $last)
$sheet3.Name = "Sheet3"
$sheet3.Activate()
